$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Update row 2 (OFICINA-CATAECSA / OTROS)
$ws.Range("D2").Value = 5596.5
$ws.Range("E2").Value = -5596.5

# Update row 4 (TOTAL)
$ws.Range("D4").Value = 18368.75
$ws.Range("E4").Value = -4645.41
$ws.Range("F4").Value = 1.338504329121045
